# Locate target paragraph (the contact-info line containing "GitHub")
$d = $word.ActiveDocument

$searchRange = $d.Content
$foundPara = $searchRange.Find.Execute("GitHub", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundPara) {
    throw "Could not locate the contact-info paragraph (GitHub marker not found)."
}

$targetPara = $searchRange.Paragraphs(1).Range
$paraStart = $targetPara.Start
$paraEnd = $targetPara.End - 1   # exclude the paragraph mark

$target = $d.Range($paraStart, $paraEnd)

$xmlChunk = @'
<?xml version="1.0" encoding="UTF-16" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
      <w:r w:rsidRPr="00B9281E">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve">Fairfax, VA | (703) 656-3648 </w:t>
      </w:r>
      <w:r w:rsidR="00B1772A">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve">| </w:t>
      </w:r>
      <w:hyperlink r:id="rId6" w:history="1">
        <w:r w:rsidR="00B1772A" w:rsidRPr="00EA2231">
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
            <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
            <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          </w:rPr>
          <w:t>ssumathi@gmu.edu</w:t>
        </w:r>
      </w:hyperlink>
      <w:r w:rsidR="00B1772A">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> | </w:t>
      </w:r>
      <w:hyperlink r:id="rId7" w:history="1">
        <w:r w:rsidR="00B1772A" w:rsidRPr="00B1772A">
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
            <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          </w:rPr>
          <w:t>LinkedIn</w:t>
        </w:r>
      </w:hyperlink>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Hyperlink"/>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          <w:u w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r w:rsidR="00B1772A">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve">| </w:t>
      </w:r>
      <w:hyperlink r:id="rId8" w:history="1">
        <w:r w:rsidR="00B1772A" w:rsidRPr="00B1772A">
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
            <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
            <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          </w:rPr>
          <w:t>GitHub</w:t>
        </w:r>
      </w:hyperlink>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Hyperlink"/>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          <w:u w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Hyperlink"/>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:u w:val="none"/>
        </w:rPr>
        <w:t>|</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="Hyperlink"/>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          <w:u w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:hyperlink r:id="rIdWebsiteLink" w:history="1">
        <w:r>
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
            <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
            <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          </w:rPr>
          <w:t>We</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
            <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
            <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          </w:rPr>
          <w:t>b</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
            <w:rFonts w:ascii="Calibri" w:eastAsia="Cambria" w:hAnsi="Calibri" w:cs="Calibri"/>
            <w:color w:val="0F9ED5" w:themeColor="accent4"/>
          </w:rPr>
          <w:t>site</w:t>
        </w:r>
      </w:hyperlink>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/customXml" Target="../customXml/item1.xml"/>
<Relationship Id="rId2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/numbering" Target="numbering.xml"/>
<Relationship Id="rId3" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/>
<Relationship Id="rId4" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/settings" Target="settings.xml"/>
<Relationship Id="rId5" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/webSettings" Target="webSettings.xml"/>
<Relationship Id="rId6" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="mailto:ssumathi@gmu.edu" TargetMode="External"/>
<Relationship Id="rId7" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.linkedin.com/in/santhosh-s-e/%20" TargetMode="External"/>
<Relationship Id="rId8" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/santhoshse99" TargetMode="External"/>
<Relationship Id="rIdWebsiteLink" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://santhoshse99.github.io" TargetMode="External"/>
<Relationship Id="rId9" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/fontTable" Target="fontTable.xml"/>
<Relationship Id="rId10" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/theme" Target="theme/theme1.xml"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xmlChunk)

Write-Output ("Updated paragraph text: " + $d.Paragraphs(2).Range.Text)
